$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column to be treated as text before writing, so that
# numeric-looking strings (e.g. "592.71") are not auto-coerced into
# numbers by the COM type-inference. We restore the original (default)
# formatting afterwards.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2 - Bitcoin
$ws.Range("D2").Value = "63.019.63"
$ws.Range("E2").Value = "  -2.40%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.116.46"
$ws.Range("E3").Value = "  -1.03%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.05%  "

# Row 5 - BNB
$ws.Range("D5").Value = "592.71"
$ws.Range("E5").Value = "  -3.06%  "

# Row 6 - Solana
$ws.Range("D6").Value = "136.21"
$ws.Range("E6").Value = "  -5.74%  "

# Row 8 - LidoStakedEther
$ws.Range("D8").Value = "3.111.92"
$ws.Range("E8").Value = "  -1.08%  "

# Row 9 - XRP
$ws.Range("E9").Value = "  -3.20%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  -3.97%  "

# Row 11 - Toncoin
$ws.Range("E11").Value = "  -4.19%  "

# Row 12 - Cardano
$ws.Range("E12").Value = "  -4.01%  "

# Row 13 - ShibaInu
$ws.Range("E13").Value = "  -3.90%  "

# Row 14 - Avalanche
$ws.Range("D14").Value = "34.12"
$ws.Range("E14").Value = "  -4.41%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "3.633.42"
$ws.Range("E15").Value = "  -0.84%  "

# Row 17 - now WrappedBTC (was WrappedEther)
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "62.977.40"
$ws.Range("E17").Value = "  -2.37%  "

# Row 18 - now WrappedEther (was WrappedBTC)
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.118.60"
$ws.Range("E18").Value = "  -0.83%  "

# Row 19 - Polkadot
$ws.Range("D19").Value = "6.71"
$ws.Range("E19").Value = "  -2.86%  "

# Row 20 - BitcoinCash
$ws.Range("D20").Value = "475.86"
$ws.Range("E20").Value = "  -0.69%  "

# Row 21 - Chainlink
$ws.Range("D21").Value = "14.19"
$ws.Range("E21").Value = "  -4.13%  "

# Row 22 - Polygon
$ws.Range("E22").Value = "  -4.19%  "

# Row 23 - Uniswap
$ws.Range("D23").Value = "7.70"
$ws.Range("E23").Value = "  -2.28%  "

# Row 24 - Litecoin
$ws.Range("D24").Value = "87.22"
$ws.Range("E24").Value = "  +1.90%  "

# Row 25 - InternetComputer(DFINITY)
$ws.Range("D25").Value = "13.01"
$ws.Range("E25").Value = "  -5.19%  "

# Row 26 - Dai
$ws.Range("E26").Value = "  +0.14%  "

# Row 27 - PancakeSwap
$ws.Range("D27").Value = "2.71"
$ws.Range("E27").Value = "  -2.86%  "

# Row 28 - NEARProtocol
$ws.Range("D28").Value = "7.15"
$ws.Range("E28").Value = "  -3.42%  "

# Row 29 - RenderToken
$ws.Range("D29").Value = "7.91"
$ws.Range("E29").Value = "  -7.99%  "

# Row 30 - ImmutableX
$ws.Range("E30").Value = "  -1.26%  "

# Row 31 - FirstDigitalUSD
$ws.Range("E31").Value = "  +0.05%  "

# Row 32 - EthereumClassic
$ws.Range("D32").Value = "26.72"
$ws.Range("E32").Value = "  -0.49%  "

# Row 33 - Hedera
$ws.Range("E33").Value = "  -8.49%  "

# Row 34 - Stacks
$ws.Range("D34").Value = "2.53"
$ws.Range("E34").Value = "  -4.57%  "

# Row 35 - Mantle
$ws.Range("E35").Value = "  -3.35%  "

# Row 36 - Filecoin
$ws.Range("E36").Value = "  -2.54%  "

# Row 37 - OKB
$ws.Range("D37").Value = "51.99"
$ws.Range("E37").Value = "  -1.60%  "

# Row 38 - PEPE
$ws.Range("D38").Value = "0.0₃0714"
$ws.Range("E38").Value = "  -5.29%  "

# Row 39 - VeChain
$ws.Range("D39").Value = "0.0388"
$ws.Range("E39").Value = "  -2.42%  "

# Row 40 - Bittensor
$ws.Range("D40").Value = "421.15"
$ws.Range("E40").Value = "  -7.93%  "

# Row 41 - Kaspa
$ws.Range("E41").Value = "  -1.41%  "

# Row 42 - Cosmos
$ws.Range("D42").Value = "8.25"
$ws.Range("E42").Value = "  -1.34%  "

# Row 43 - dogwifhat
$ws.Range("E43").Value = "  -11.76%  "

# Row 44 - Maker
$ws.Range("D44").Value = "2.884.17"
$ws.Range("E44").Value = "  +0.06%  "

# Row 45 - TheGraph
$ws.Range("E45").Value = "  +0.44%  "

# Row 46 - Fetch.AI
$ws.Range("E46").Value = "  -5.85%  "

# Row 47 - USDe
$ws.Range("D47").Value = "0.999"
$ws.Range("E47").Value = "  -0.06%  "

# Row 48 - InjectiveProtocol
$ws.Range("D48").Value = "25.72"
$ws.Range("E48").Value = "  -3.47%  "

# Row 49 - ThetaToken
$ws.Range("E49").Value = "  -7.10%  "

# Row 50 - Stellar
$ws.Range("E50").Value = "  -1.22%  "

# Row 51 - Monero
$ws.Range("D51").Value = "119.62"
$ws.Range("E51").Value = "  -1.46%  "

# Restore default formatting on the Price column now that the literal
# text values are committed.
$ws.Range("D2:D51").ClearFormats()
